# Add a header row (Email / Password) above the existing login data on the
# "Login" sheet, center the new headers, capitalize the password value, and
# keep the existing hyperlinks pointing at the (now shifted) data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Push the existing data down one row so a header row can be inserted at the top.
$ws.Rows.Item(1).Insert()

# Populate the new header row.
$ws.Range("A1").Value2 = "Email"
$ws.Range("B1").Value2 = "Password"

# The password text changes case (saku@123 -> Saku@123); it now lives in B2.
$ws.Range("B2").Value2 = "Saku@123"

# Center-align the new header cells.
$ws.Range("A1:B1").HorizontalAlignment = -4108   # xlCenter

# The row-insert does not relocate the worksheet's hyperlink anchors, so they
# are still pointing at row 1. Remove the stale anchors and recreate them on
# the data row (row 2), preserving the original target addresses.
$ws.Range("B1").Hyperlinks.Delete()
$ws.Range("A1").Hyperlinks.Delete()

# Reset to the default style first so both new hyperlinks are derived from the
# same base format (keeps the stylesheet as small as possible), then recreate
# the links and restore the normal "Hyperlink" look.
$ws.Range("A2:B2").Style = "Normal"
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:Saku@123") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sakunthalanfm@gmail.com") | Out-Null
$ws.Range("A2:B2").Style = "Hyperlink"

# Update the saved selection.
$ws.Range("A8").Select() | Out-Null

Write-Host "Edit complete"
